# Updates the match/odds data for the "Poland Ekstraklasa" league sheet.
# Rows 147,148,149,150,151,152,154 (by sheet row number) get their match
# data (id in col B, HomeTeam/AwayTeam/FTHG/FTAG/FTR in cols F:J, and all
# the odds columns K:AC) reassigned to match an updated upstream data
# snapshot. Column A (running index), C/D (Div / Div Original Name) and
# E (Date) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland Ekstraklasa")

function Set-MatchRow {
    param($Row, $Id, $HomeTeam, $AwayTeam, $FTHG, $FTAG, $FTR, $Odds)

    $ws.Cells.Item($Row, 2).Value = $Id          # B - id
    $ws.Cells.Item($Row, 6).Value = $HomeTeam    # F - HomeTeam
    $ws.Cells.Item($Row, 7).Value = $AwayTeam    # G - AwayTeam
    $ws.Cells.Item($Row, 8).Value = $FTHG        # H - FTHG
    $ws.Cells.Item($Row, 9).Value = $FTAG        # I - FTAG
    $ws.Cells.Item($Row, 10).Value = $FTR        # J - FTR

    # K..AC (19 odds columns -> column index 11..29)
    for ($i = 0; $i -lt $Odds.Length; $i++) {
        $ws.Cells.Item($Row, 11 + $i).Value = $Odds[$i]
    }
}

# Row 147 <- old row 150 data
Set-MatchRow 147 5460884 "Miedz Legnica" "Gornik Zabrze" 0 0 "D" `
    @(3.6, 3.5, 1.909, 3, 3.5, 2.1, 0.25, 1.95, 1.9, 2.75, 1.975, 1.875, -1, 2.5, -1, 0.475, -0.5, -1, 0.875)

# Row 148 <- old row 151 data
Set-MatchRow 148 5456603 "Lech Poznan" "Jagiellonia Bialystok" 2 0 "H" `
    @(1.363, 4.75, 6.5, 1.222, 5.5, 8, -1.75, 1.925, 1.925, 3.25, 1.95, 1.9, 0.222, -1, -1, 0.4625, -0.5, -1, 0.8999999999999999)

# Row 149 <- old row 154 data
Set-MatchRow 149 5456594 "Rakow Czestochowa" "Zaglebie Lubin" 1 1 "D" `
    @(1.444, 4.5, 5.75, 1.3, 5.25, 7, -1.5, 1.9, 1.95, 3, 1.9, 1.95, -1, 4.25, -1, -1, 0.95, -1, 0.95)

# Row 150 <- old row 152 data
Set-MatchRow 150 5428774 "Pogon Szczecin" "Radomiak Radom" 4 0 "H" `
    @(1.571, 4, 4.75, 1.533, 4.333, 4.75, -1, 1.875, 1.975, 3, 1.875, 1.975, 0.5329999999999999, -1, -1, 0.875, -1, 0.875, -1)

# Row 151 <- old row 147 data
Set-MatchRow 151 5461475 "Widzew Lodz" "Korona Kielce" 0 3 "A" `
    @(2.1, 3.3, 3.2, 2.375, 3.3, 2.7, 0, 1.8, 2.05, 2.5, 1.825, 2.025, -1, -1, 1.7, -1, 1.05, 0.825, -1)

# Row 152 <- old row 148 data
Set-MatchRow 152 5461474 "Legia Warsaw" "Slask Wroclaw" 3 1 "H" `
    @(1.7, 3.8, 4, 1.833, 3.8, 3.4, -0.5, 1.825, 2.025, 2.75, 1.9, 1.95, 0.833, -1, -1, 0.825, -1, 0.8999999999999999, -1)

# Row 154 <- old row 149 data
Set-MatchRow 154 5467427 "Stal Mielec" "Warta Poznan" 1 0 "H" `
    @(2.375, 3.2, 2.8, 2.6, 3.1, 2.625, 0, 1.925, 1.925, 2.25, 1.975, 1.875, 1.6, -1, -1, 0.925, -1, -1, 0.875)
